$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 1.95
$ws.Range("J2").Value = 4.5
$ws.Range("L2").Value = 2.63
$ws.Range("X2").Value = 19
$ws.Range("Y2").Value = 13
$ws.Range("AA2").Value = 34
$ws.Range("AI2").Value = 8.5
$ws.Range("AK2").Value = 17
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 251
$ws.Range("AY2").Value = 11
$ws.Range("BA2").Value = 41
$ws.Range("BB2").Value = 67
